$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $r = $ws.Range($rangeAddr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '59.685.61'
Set-TextValue 'D3' '2.587.69'
Set-TextValue 'E3' '  +0.71%  '
Set-TextValue 'E4' '  +0.07%  '
Set-TextValue 'D5' '556.41'
Set-TextValue 'E5' '  -0.96%  '
Set-TextValue 'D6' '141.06'
Set-TextValue 'E6' '  -1.10%  '
Set-TextValue 'E7' '  -0.12%  '
Set-TextValue 'D8' '0.596'
Set-TextValue 'E8' '  -0.16%  '
Set-TextValue 'D9' '2.605.76'
Set-TextValue 'E9' '  +1.20%  '
Set-TextValue 'E10' '  +0.89%  '
Set-TextValue 'E11' '  +1.59%  '
Set-TextValue 'E12' '  +6.83%  '
Set-TextValue 'D13' '0.360'
Set-TextValue 'E13' '  +5.72%  '
Set-TextValue 'D14' '3.049.99'
Set-TextValue 'E14' '  +0.85%  '
Set-TextValue 'B15' 'WrappedBTC'
Set-TextValue 'C15' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D15' '59.650.05'
Set-TextValue 'E15' '  +1.35%  '
Set-TextValue 'B16' 'Avalanche'
Set-TextValue 'C16' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D16' '23.27'
Set-TextValue 'E16' '  +6.31%  '
Set-TextValue 'E17' '  +0.57%  '
Set-TextValue 'D18' '2.599.96'
Set-TextValue 'E18' '  +0.74%  '
Set-TextValue 'D19' '4.57'
Set-TextValue 'E19' '  +1.64%  '
Set-TextValue 'D20' '341.16'
Set-TextValue 'E20' '  +1.93%  '
Set-TextValue 'D21' '10.51'
Set-TextValue 'E21' '  +3.68%  '
Set-TextValue 'D22' '6.69'
Set-TextValue 'E22' '  +8.98%  '
Set-TextValue 'D23' '0.999'
Set-TextValue 'E23' '  +0.00%  '
Set-TextValue 'D24' '0.497'
Set-TextValue 'E24' '  +11.11%  '
Set-TextValue 'D25' '62.46'
Set-TextValue 'E25' '  -1.73%  '
Set-TextValue 'D26' '0.995'
Set-TextValue 'E26' '  -0.64%  '
Set-TextValue 'D27' '0.160'
Set-TextValue 'E27' '  -0.30%  '
Set-TextValue 'D28' '7.52'
Set-TextValue 'E28' '  +4.24%  '
Set-TextValue 'D29' '0.0₃0776'
Set-TextValue 'E29' '  -0.26%  '
Set-TextValue 'E30' '  -0.07%  '
Set-TextValue 'D31' '1.70'
Set-TextValue 'E31' '  +1.54%  '
Set-TextValue 'E32' '  +2.01%  '
Set-TextValue 'D33' '158.66'
Set-TextValue 'E33' '  +0.15%  '
Set-TextValue 'D34' '19.34'
Set-TextValue 'E34' '  +2.15%  '
Set-TextValue 'D35' '4.08'
Set-TextValue 'E35' '  +2.56%  '
Set-TextValue 'D36' '0.913'
Set-TextValue 'E36' '  +4.34%  '
Set-TextValue 'E37' '  +4.39%  '
Set-TextValue 'D38' '37.76'
Set-TextValue 'E38' '  +2.63%  '
Set-TextValue 'D39' '1.52'
Set-TextValue 'E39' '  +2.46%  '
Set-TextValue 'D40' '0.843'
Set-TextValue 'E40' '  -3.57%  '
Set-TextValue 'E41' '  +1.83%  '
Set-TextValue 'D42' '291.95'
Set-TextValue 'E42' '  +0.31%  '
Set-TextValue 'D43' '136.55'
Set-TextValue 'E43' '  +10.36%  '
Set-TextValue 'D44' '0.998'
Set-TextValue 'E44' '  -0.16%  '
Set-TextValue 'D45' '0.0977'
Set-TextValue 'E45' '  +0.72%  '
Set-TextValue 'E46' '  +1.20%  '
Set-TextValue 'E47' '  +3.31%  '
Set-TextValue 'D48' '0.0537'
Set-TextValue 'E48' '  +1.19%  '
Set-TextValue 'D49' '10.65'
Set-TextValue 'E49' '  +0.41%  '
Set-TextValue 'E50' '  +7.41%  '
Set-TextValue 'D51' '18.84'
Set-TextValue 'E51' '  +2.30%  '
